$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column D (shifts old D:K -> F:M)
$ws.Range("D:E").Insert()

# Copy number formats/styles from column F (the old D, now shifted) into new D:E
# so the new columns inherit the same per-row styles (date format row 7/38/80,
# number format elsewhere) instead of a default style.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)

# Populate the two new columns (newest quarters: 2018-12-31 and 2018-09-30)
# with the reported financial figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 301400
$ws.Range("E8").Value = 242800
$ws.Range("D9").Value = 73100
$ws.Range("E9").Value = 55100
$ws.Range("D10").Value = 228300
$ws.Range("E10").Value = 187700
$ws.Range("D12").Value = 54600
$ws.Range("E12").Value = 148200
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 300
$ws.Range("D15").Value = 121200
$ws.Range("E15").Value = 80000
$ws.Range("D17").Value = 270300
$ws.Range("E17").Value = 308900
$ws.Range("D18").Value = 31100
$ws.Range("E18").Value = -66100
$ws.Range("D20").Value = 288200
$ws.Range("E20").Value = -26000
$ws.Range("D21").Value = 442800
$ws.Range("E21").Value = -9700
$ws.Range("D22").Value = 32200
$ws.Range("E22").Value = 22600
$ws.Range("D23").Value = 287000
$ws.Range("E23").Value = -114700
$ws.Range("D24").Value = 101500
$ws.Range("E24").Value = 11400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 185600
$ws.Range("E26").Value = -126100
$ws.Range("D27").Value = 185600
$ws.Range("E27").Value = -126100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -288200
$ws.Range("E32").Value = 26000
$ws.Range("D33").Value = 185600
$ws.Range("E33").Value = -126100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 185600
$ws.Range("E35").Value = -126100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 173500
$ws.Range("E41").Value = 192600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 140000
$ws.Range("E43").Value = 252300
$ws.Range("D44").Value = 84800
$ws.Range("E44").Value = 90000
$ws.Range("D45").Value = 111400
$ws.Range("E45").Value = 105800
$ws.Range("D46").Value = 509700
$ws.Range("E46").Value = 640800
$ws.Range("D47").Value = 66200
$ws.Range("E47").Value = 103100
$ws.Range("D48").Value = 3459700
$ws.Range("E48").Value = 3509500
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 52600
$ws.Range("E52").Value = 76000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 4088200
$ws.Range("E54").Value = 4329500
$ws.Range("D57").Value = 176500
$ws.Range("E57").Value = 153900
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 207800
$ws.Range("E59").Value = 474500
$ws.Range("D60").Value = 384300
$ws.Range("E60").Value = 628400
$ws.Range("D61").Value = 2120500
$ws.Range("E61").Value = 2094500
$ws.Range("D62").Value = 641900
$ws.Range("E62").Value = 671500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 3146700
$ws.Range("E66").Value = 3394500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -1167200
$ws.Range("E72").Value = -1352800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 941500
$ws.Range("E76").Value = 934900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 185600
$ws.Range("E81").Value = -126100
$ws.Range("D83").Value = 123500
$ws.Range("E83").Value = 82400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 170200
$ws.Range("E89").Value = 89800
$ws.Range("D91").Value = -4400
$ws.Range("E91").Value = -700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -26800
$ws.Range("E94").Value = -942800
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -165300
$ws.Range("E100").Value = 914000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -21900
$ws.Range("E102").Value = 61000

# Row 91 (Capital Expenditures) was also restated for several of the
# older quarters in addition to gaining the two new columns.
$ws.Range("F91").Value = -1100
$ws.Range("G91").Value = -1800
$ws.Range("H91").Value = -1200
$ws.Range("I91").Value = -200
$ws.Range("J91").Value = -1200
